# Add team record (Wins/Losses/Ties) columns to the sheet, matching the
# diff that introduces AD/AE/AF columns:
#   AD1="Wins", AE1="Losses", AF1="Ties" (same bold header style as A1:AC1)
#   AD2:AD58 = 87, AE2:AE58 = 75, AF2:AF58 = 0 (numeric)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header values for the new columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the existing header row (A1:AC1),
# which carries the bold/centered style (style index 1 in the original file).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the team record for every player row (2 through 58).
for ($r = 2; $r -le 58; $r++) {
    $ws.Cells.Item($r, 30).Value = 87   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 75   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
